$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: extend hours and grow the task note text ---
$ws.Range("B35").Value = 4.5
$ws.Range("D35").Value = "Verified no PW's in unexpected places`nWeek 6: submitted professional dev feedback, installed Putty, worked through the rest of AWS exercise videos but left off on actual exercise at the point of getting database set up (step 5).  Hoping to do step 6 after more project work is complete."
$ws.Rows(35).RowHeight = 45

# --- Row 36: was an empty placeholder date row, now a real logged day ---
$ws.Range("A36").Value = 43530
$ws.Range("A36").NumberFormat = "d-mmm"
$ws.Range("B36").Value = 3.5
$ws.Range("D36").Value = "Indie Project: organized priorities, figured out how to generate UML diagrams from IntelliJ, tried and failed to recreate Log4J problem, set up Servlet shells`nWeek 7: Intro video"
$ws.Rows(36).RowHeight = 45

# --- Insert four fresh blank rows right after row 37 (pushing the
#     "Issues/Loose Ends" block and everything below it down by 4),
#     anchored so the new rows don't inherit row 37's formatting ---
$ws.Rows("39:42").Insert()

# Row 37 gains a second (still empty) cell carrying the date/wrap style
$ws.Range("D37").NumberFormat = "h:mm"
$ws.Range("D37").WrapText = $true

# Restore the active sheet view/selection to match the edited document
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("D38").Select()
